$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 data (appended entry).
# A4 and C4 look like numbers ("222" / "2000") but must stay text, like the
# rest of the sheet's numeric-looking values (e.g. C2="12", C3="2323"), so
# force a text number format on those two cells before writing the values.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("C4").NumberFormat = "@"

$ws.Range("A4").Value = "222"
$ws.Range("B4").Value = "احمد"
$ws.Range("C4").Value = "2000"
$ws.Range("D4").Value = "الجزائري"
$ws.Range("E4").Value = "الرحلة 1"
$ws.Range("F4").Value = "C3"
$ws.Range("G4").Value = "NRC"
$ws.Range("H4").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٤:٢١:٢٢ م"
